# Update the dSF column (F) values to reflect the repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 5
    5  = 0
    8  = 1
    11 = 0
    13 = 0
    16 = 2
    18 = 1
    20 = -5
    21 = -4
    22 = 4
    23 = -4
    24 = 1
    28 = -3
    29 = -2
    32 = 0
    33 = 3
    35 = 1
    36 = 2
    37 = -1
    38 = -7
    40 = 0
    43 = -6
    45 = 1
    47 = 0
    50 = 1
    51 = -7
    54 = -2
    57 = 6
    58 = 0
    59 = -6
    63 = -8
    65 = -2
    69 = -5
    70 = -2
    72 = -5
    74 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
